# Remove the "82k resistor 1%" (CRGCQ1206F82K / A129833CT-ND) line item from
# row 5 of every vendor sheet (Digikey, Arrow, Mouser). The new schematic no
# longer uses this part, so its row is wiped back to a blank spacer row (same
# as rows 4/10/etc elsewhere on these sheets) rather than being physically
# deleted - formatting/styles stay, values/formula/hyperlink go.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Digikey", "Arrow", "Mouser")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Wipe the row-5 line item back to an empty spacer row, keeping cell
    # formatting (ClearContents leaves styles alone, unlike Clear()).
    $ws.Range("A5:K5").ClearContents()

    # That row's A5 cell carried a hyperlink to the now-removed part's
    # product page - drop it now that A5 is blank. Find it first, then
    # delete outside the loop so we don't mutate the collection while
    # iterating it.
    $hyperlinkToRemove = $null
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq '$A$5') {
            $hyperlinkToRemove = $hl
        }
    }
    if ($hyperlinkToRemove -ne $null) {
        $hyperlinkToRemove.Delete()
    }
}

# Restore each sheet's former selection rectangle (A5:K5, the row that just
# got cleared) ...
$wsDigikey = $wb.Worksheets.Item("Digikey")
[void]$wsDigikey.Range("A5:K5").Select()

$wsArrow = $wb.Worksheets.Item("Arrow")
[void]$wsArrow.Range("A5:K5").Select()

# ... except Mouser, which becomes the active tab with a different selection.
$wsMouser = $wb.Worksheets.Item("Mouser")
$wsMouser.Activate()
[void]$wsMouser.Range("N15").Select()
